$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 156
$ws.Cells.Item(156, 2).Value = 1
$c = $ws.Cells.Item(156, 3)
$c.NumberFormat = "@"
$c.Value = "5"
$c.NumberFormat = "General"
$ws.Cells.Item(156, 8).Value = 0.02583063914096329
$ws.Cells.Item(156, 9).Value = 0.003612079715285082
$ws.Cells.Item(156, 10).Value = 2.583063914096329
$ws.Cells.Item(156, 11).Value = 0.3612079715285083
$ws.Cells.Item(156, 12).Value = 7.151181916516648

# Row 157
$ws.Cells.Item(157, 2).Value = 1
$c = $ws.Cells.Item(157, 3)
$c.NumberFormat = "@"
$c.Value = "8"
$c.NumberFormat = "General"
$ws.Cells.Item(157, 8).Value = 0.03042928712297144
$ws.Cells.Item(157, 9).Value = 0.002026808196544004
$ws.Cells.Item(157, 10).Value = 3.042928712297144
$ws.Cells.Item(157, 11).Value = 0.2026808196544004
$ws.Cells.Item(157, 12).Value = 15.01340243978571

# Row 158
$ws.Cells.Item(158, 2).Value = 1
$c = $ws.Cells.Item(158, 3)
$c.NumberFormat = "@"
$c.Value = "10"
$c.NumberFormat = "General"
$ws.Cells.Item(158, 8).Value = 0.04371646809885044
$ws.Cells.Item(158, 9).Value = 0.001360777663753391
$ws.Cells.Item(158, 10).Value = 4.371646809885044
$ws.Cells.Item(158, 11).Value = 0.1360777663753391
$ws.Cells.Item(158, 12).Value = 32.12609176598965

# Row 159
$ws.Cells.Item(159, 2).Value = 1
$c = $ws.Cells.Item(159, 3)
$c.NumberFormat = "@"
$c.Value = "12"
$c.NumberFormat = "General"
$ws.Cells.Item(159, 8).Value = 0.04833598122596161
$ws.Cells.Item(159, 9).Value = 0.001487423355882465
$ws.Cells.Item(159, 10).Value = 4.833598122596161
$ws.Cells.Item(159, 11).Value = 0.1487423355882465
$ws.Cells.Item(159, 12).Value = 32.49645168929368

# Row 160
$ws.Cells.Item(160, 2).Value = 1
$c = $ws.Cells.Item(160, 3)
$c.NumberFormat = "@"
$c.Value = "14"
$c.NumberFormat = "General"
$ws.Cells.Item(160, 8).Value = 0.05334177976485388
$ws.Cells.Item(160, 9).Value = 0.001090415921081629
$ws.Cells.Item(160, 10).Value = 5.334177976485388
$ws.Cells.Item(160, 11).Value = 0.1090415921081629
$ws.Cells.Item(160, 12).Value = 48.91874626329918

# Row 161
$ws.Cells.Item(161, 2).Value = 1
$c = $ws.Cells.Item(161, 3)
$c.NumberFormat = "@"
$c.Value = "16"
$c.NumberFormat = "General"
$ws.Cells.Item(161, 8).Value = 0.05543920456722629
$ws.Cells.Item(161, 9).Value = 0.0008346471988880426
$ws.Cells.Item(161, 10).Value = 5.543920456722629
$ws.Cells.Item(161, 11).Value = 0.08346471988880426
$ws.Cells.Item(161, 12).Value = 66.4223214803633

# Row 162
$ws.Cells.Item(162, 2).Value = 1
$c = $ws.Cells.Item(162, 3)
$c.NumberFormat = "@"
$c.Value = "5"
$c.NumberFormat = "General"
$ws.Cells.Item(162, 8).Value = 0.06059726576506375
$ws.Cells.Item(162, 9).Value = 0.02054808413553775
$ws.Cells.Item(162, 10).Value = 6.059726576506375
$ws.Cells.Item(162, 11).Value = 2.054808413553775
$ws.Cells.Item(162, 12).Value = 2.94904699461792

# Row 163
$ws.Cells.Item(163, 2).Value = 1
$c = $ws.Cells.Item(163, 3)
$c.NumberFormat = "@"
$c.Value = "8"
$c.NumberFormat = "General"
$ws.Cells.Item(163, 8).Value = 0.087618499403884
$ws.Cells.Item(163, 9).Value = 0.009720798641250241
$ws.Cells.Item(163, 10).Value = 8.7618499403884
$ws.Cells.Item(163, 11).Value = 0.9720798641250241
$ws.Cells.Item(163, 12).Value = 9.013508317318148

# Row 164
$ws.Cells.Item(164, 2).Value = 1
$c = $ws.Cells.Item(164, 3)
$c.NumberFormat = "@"
$c.Value = "10"
$c.NumberFormat = "General"
$ws.Cells.Item(164, 8).Value = 0.09617722100753023
$ws.Cells.Item(164, 9).Value = 0.00525114836410793
$ws.Cells.Item(164, 10).Value = 9.617722100753024
$ws.Cells.Item(164, 11).Value = 0.525114836410793
$ws.Cells.Item(164, 12).Value = 18.3154644162999

# Row 165
$ws.Cells.Item(165, 2).Value = 1
$c = $ws.Cells.Item(165, 3)
$c.NumberFormat = "@"
$c.Value = "12"
$c.NumberFormat = "General"
$ws.Cells.Item(165, 8).Value = 0.09097484823570245
$ws.Cells.Item(165, 9).Value = 0.003837407134138946
$ws.Cells.Item(165, 10).Value = 9.097484823570245
$ws.Cells.Item(165, 11).Value = 0.3837407134138946
$ws.Cells.Item(165, 12).Value = 23.70737455151883

# Row 166
$ws.Cells.Item(166, 2).Value = 1
$c = $ws.Cells.Item(166, 3)
$c.NumberFormat = "@"
$c.Value = "14"
$c.NumberFormat = "General"
$ws.Cells.Item(166, 8).Value = 0.08888059119759428
$ws.Cells.Item(166, 9).Value = 0.002597295191077616
$ws.Cells.Item(166, 10).Value = 8.88805911975943
$ws.Cells.Item(166, 11).Value = 0.2597295191077617
$ws.Cells.Item(166, 12).Value = 34.2204426754888

# Row 167
$ws.Cells.Item(167, 2).Value = 1
$c = $ws.Cells.Item(167, 3)
$c.NumberFormat = "@"
$c.Value = "16"
$c.NumberFormat = "General"
$ws.Cells.Item(167, 8).Value = 0.08894038457169606
$ws.Cells.Item(167, 9).Value = 0.001915186785349299
$ws.Cells.Item(167, 10).Value = 8.894038457169607
$ws.Cells.Item(167, 11).Value = 0.1915186785349299
$ws.Cells.Item(167, 12).Value = 46.4395354291643

# Row 168
$ws.Cells.Item(168, 2).Value = 1
$c = $ws.Cells.Item(168, 3)
$c.NumberFormat = "@"
$c.Value = "5"
$c.NumberFormat = "General"
$ws.Cells.Item(168, 8).Value = 0.02816781695331816
$ws.Cells.Item(168, 9).Value = 0.003579348361891447
$ws.Cells.Item(168, 10).Value = 2.816781695331816
$ws.Cells.Item(168, 11).Value = 0.3579348361891447
$ws.Cells.Item(168, 12).Value = 7.869537721786136

# Row 169
$ws.Cells.Item(169, 2).Value = 1
$c = $ws.Cells.Item(169, 3)
$c.NumberFormat = "@"
$c.Value = "8"
$c.NumberFormat = "General"
$ws.Cells.Item(169, 8).Value = 0.03164037831480759
$ws.Cells.Item(169, 9).Value = 0.002010211181439609
$ws.Cells.Item(169, 10).Value = 3.164037831480759
$ws.Cells.Item(169, 11).Value = 0.201021118143961
$ws.Cells.Item(169, 12).Value = 15.7398280374445

# Row 170
$ws.Cells.Item(170, 2).Value = 1
$c = $ws.Cells.Item(170, 3)
$c.NumberFormat = "@"
$c.Value = "10"
$c.NumberFormat = "General"
$ws.Cells.Item(170, 8).Value = 0.04364840924816327
$ws.Cells.Item(170, 9).Value = 0.001361576528699496
$ws.Cells.Item(170, 10).Value = 4.364840924816327
$ws.Cells.Item(170, 11).Value = 0.1361576528699496
$ws.Cells.Item(170, 12).Value = 32.05725739841731

# Row 171
$ws.Cells.Item(171, 2).Value = 1
$c = $ws.Cells.Item(171, 3)
$c.NumberFormat = "@"
$c.Value = "12"
$c.NumberFormat = "General"
$ws.Cells.Item(171, 8).Value = 0.04683682271467493
$ws.Cells.Item(171, 9).Value = 0.001511023127474993
$ws.Cells.Item(171, 10).Value = 4.683682271467493
$ws.Cells.Item(171, 11).Value = 0.1511023127474993
$ws.Cells.Item(171, 12).Value = 30.99676097806787

# Row 172
$ws.Cells.Item(172, 2).Value = 1
$c = $ws.Cells.Item(172, 3)
$c.NumberFormat = "@"
$c.Value = "14"
$c.NumberFormat = "General"
$ws.Cells.Item(172, 8).Value = 0.05041818775781315
$ws.Cells.Item(172, 9).Value = 0.001130535467917572
$ws.Cells.Item(172, 10).Value = 5.041818775781315
$ws.Cells.Item(172, 11).Value = 0.1130535467917572
$ws.Cells.Item(172, 12).Value = 44.59673242333798

# Row 173
$ws.Cells.Item(173, 2).Value = 1
$c = $ws.Cells.Item(173, 3)
$c.NumberFormat = "@"
$c.Value = "16"
$c.NumberFormat = "General"
$ws.Cells.Item(173, 8).Value = 0.05187252026739086
$ws.Cells.Item(173, 9).Value = 0.0008781215385135165
$ws.Cells.Item(173, 10).Value = 5.187252026739086
$ws.Cells.Item(173, 11).Value = 0.08781215385135165
$ws.Cells.Item(173, 12).Value = 59.07214205814907

# Row 174
$ws.Cells.Item(174, 2).Value = 1
$c = $ws.Cells.Item(174, 3)
$c.NumberFormat = "@"
$c.Value = "5"
$c.NumberFormat = "General"
$ws.Cells.Item(174, 8).Value = 0.06075528085934945
$ws.Cells.Item(174, 9).Value = 0.01895616287462724
$ws.Cells.Item(174, 10).Value = 6.075528085934945
$ws.Cells.Item(174, 11).Value = 1.895616287462724
$ws.Cells.Item(174, 12).Value = 3.205041086699576

# Row 175
$ws.Cells.Item(175, 2).Value = 1
$c = $ws.Cells.Item(175, 3)
$c.NumberFormat = "@"
$c.Value = "8"
$c.NumberFormat = "General"
$ws.Cells.Item(175, 8).Value = 0.087618499403884
$ws.Cells.Item(175, 9).Value = 0.009720798641250241
$ws.Cells.Item(175, 10).Value = 8.7618499403884
$ws.Cells.Item(175, 11).Value = 0.9720798641250241
$ws.Cells.Item(175, 12).Value = 9.013508317318148

# Row 176
$ws.Cells.Item(176, 2).Value = 1
$c = $ws.Cells.Item(176, 3)
$c.NumberFormat = "@"
$c.Value = "10"
$c.NumberFormat = "General"
$ws.Cells.Item(176, 8).Value = 0.09613345357484082
$ws.Cells.Item(176, 9).Value = 0.004377529766672688
$ws.Cells.Item(176, 10).Value = 9.613345357484082
$ws.Cells.Item(176, 11).Value = 0.4377529766672688
$ws.Cells.Item(176, 12).Value = 21.96066245093995

# Row 177
$ws.Cells.Item(177, 2).Value = 1
$c = $ws.Cells.Item(177, 3)
$c.NumberFormat = "@"
$c.Value = "12"
$c.NumberFormat = "General"
$ws.Cells.Item(177, 8).Value = 0.0909428646863728
$ws.Cells.Item(177, 9).Value = 0.004798306055094884
$ws.Cells.Item(177, 10).Value = 9.09428646863728
$ws.Cells.Item(177, 11).Value = 0.4798306055094884
$ws.Cells.Item(177, 12).Value = 18.95311879695728

# Row 178
$ws.Cells.Item(178, 2).Value = 1
$c = $ws.Cells.Item(178, 3)
$c.NumberFormat = "@"
$c.Value = "14"
$c.NumberFormat = "General"
$ws.Cells.Item(178, 8).Value = 0.08885697609442644
$ws.Cells.Item(178, 9).Value = 0.004251317853824058
$ws.Cells.Item(178, 10).Value = 8.885697609442644
$ws.Cells.Item(178, 11).Value = 0.4251317853824059
$ws.Cells.Item(178, 12).Value = 20.9010427236109

# Row 179
$ws.Cells.Item(179, 2).Value = 1
$c = $ws.Cells.Item(179, 3)
$c.NumberFormat = "@"
$c.Value = "16"
$c.NumberFormat = "General"
$ws.Cells.Item(179, 8).Value = 0.08892297169456853
$ws.Cells.Item(179, 9).Value = 0.004876190379901471
$ws.Cells.Item(179, 10).Value = 8.892297169456853
$ws.Cells.Item(179, 11).Value = 0.4876190379901471
$ws.Cells.Item(179, 12).Value = 18.23615666465535

# Row 249
$ws.Cells.Item(249, 2).Value = 1
$c = $ws.Cells.Item(249, 3)
$c.NumberFormat = "@"
$c.Value = "11"
$c.NumberFormat = "General"
$ws.Cells.Item(249, 8).Value = 0.04885867158874491
$ws.Cells.Item(249, 9).Value = 0.006712104054796181
$ws.Cells.Item(249, 10).Value = 4.885867158874491
$ws.Cells.Item(249, 11).Value = 0.6712104054796182
$ws.Cells.Item(249, 12).Value = 7.279188640383576

# Row 250
$ws.Cells.Item(250, 2).Value = 1
$c = $ws.Cells.Item(250, 3)
$c.NumberFormat = "@"
$c.Value = "14"
$c.NumberFormat = "General"
$ws.Cells.Item(250, 8).Value = 0.06421266112280444
$ws.Cells.Item(250, 9).Value = 0.004292201089144836
$ws.Cells.Item(250, 10).Value = 6.421266112280444
$ws.Cells.Item(250, 11).Value = 0.4292201089144836
$ws.Cells.Item(250, 12).Value = 14.96031052347502

# Row 252
$ws.Cells.Item(252, 2).Value = 1
$c = $ws.Cells.Item(252, 3)
$c.NumberFormat = "@"
$c.Value = "11"
$c.NumberFormat = "General"
$ws.Cells.Item(252, 8).Value = 0.05219047997422743
$ws.Cells.Item(252, 9).Value = 0.004969022337540626
$ws.Cells.Item(252, 10).Value = 5.219047997422743
$ws.Cells.Item(252, 11).Value = 0.4969022337540626
$ws.Cells.Item(252, 12).Value = 10.50316871790491

# Row 253
$ws.Cells.Item(253, 2).Value = 1
$c = $ws.Cells.Item(253, 3)
$c.NumberFormat = "@"
$c.Value = "14"
$c.NumberFormat = "General"
$ws.Cells.Item(253, 8).Value = 0.074089907926705
$ws.Cells.Item(253, 9).Value = 0.008278936035658416
$ws.Cells.Item(253, 10).Value = 7.4089907926705
$ws.Cells.Item(253, 11).Value = 0.8278936035658416
$ws.Cells.Item(253, 12).Value = 8.949206468994381

# Row 255
$ws.Cells.Item(255, 2).Value = 1
$c = $ws.Cells.Item(255, 3)
$c.NumberFormat = "@"
$c.Value = "11"
$c.NumberFormat = "General"
$ws.Cells.Item(255, 8).Value = 0.04365356767031847
$ws.Cells.Item(255, 9).Value = 0.01147712867672637
$ws.Cells.Item(255, 10).Value = 4.365356767031847
$ws.Cells.Item(255, 11).Value = 1.147712867672637
$ws.Cells.Item(255, 12).Value = 3.803526901187432

# Row 256
$ws.Cells.Item(256, 2).Value = 1
$c = $ws.Cells.Item(256, 3)
$c.NumberFormat = "@"
$c.Value = "14"
$c.NumberFormat = "General"
$ws.Cells.Item(256, 8).Value = 0.08831198262449647
$ws.Cells.Item(256, 9).Value = 0.0128792168863935
$ws.Cells.Item(256, 10).Value = 8.831198262449647
$ws.Cells.Item(256, 11).Value = 1.28792168863935
$ws.Cells.Item(256, 12).Value = 6.856937297002535

# Row 258
$ws.Cells.Item(258, 2).Value = 1
$c = $ws.Cells.Item(258, 3)
$c.NumberFormat = "@"
$c.Value = "11"
$c.NumberFormat = "General"
$ws.Cells.Item(258, 8).Value = 0.06600493542018016
$ws.Cells.Item(258, 9).Value = 0.0131191742594598
$ws.Cells.Item(258, 10).Value = 6.600493542018016
$ws.Cells.Item(258, 11).Value = 1.31191742594598
$ws.Cells.Item(258, 12).Value = 5.031180630334733

# Row 259
$ws.Cells.Item(259, 2).Value = 1
$c = $ws.Cells.Item(259, 3)
$c.NumberFormat = "@"
$c.Value = "14"
$c.NumberFormat = "General"
$ws.Cells.Item(259, 8).Value = 0.07767775539242883
$ws.Cells.Item(259, 9).Value = 0.00804629870762249
$ws.Cells.Item(259, 10).Value = 7.767775539242883
$ws.Cells.Item(259, 11).Value = 0.804629870762249
$ws.Cells.Item(259, 12).Value = 9.653849330604947

# Row 261
$ws.Cells.Item(261, 2).Value = 1
$c = $ws.Cells.Item(261, 3)
$c.NumberFormat = "@"
$c.Value = "11"
$c.NumberFormat = "General"
$ws.Cells.Item(261, 8).Value = 0.05699844292981404
$ws.Cells.Item(261, 9).Value = 0.007041336165540353
$ws.Cells.Item(261, 10).Value = 5.699844292981404
$ws.Cells.Item(261, 11).Value = 0.7041336165540353
$ws.Cells.Item(261, 12).Value = 8.094833365400042

# Row 262
$ws.Cells.Item(262, 2).Value = 1
$c = $ws.Cells.Item(262, 3)
$c.NumberFormat = "@"
$c.Value = "14"
$c.NumberFormat = "General"
$ws.Cells.Item(262, 8).Value = 0.08878609080926925
$ws.Cells.Item(262, 9).Value = 0.009006989176214581
$ws.Cells.Item(262, 10).Value = 8.878609080926925
$ws.Cells.Item(262, 11).Value = 0.9006989176214582
$ws.Cells.Item(262, 12).Value = 9.857466137933551

# Row 264
$ws.Cells.Item(264, 2).Value = 1
$c = $ws.Cells.Item(264, 3)
$c.NumberFormat = "@"
$c.Value = "11"
$c.NumberFormat = "General"
$ws.Cells.Item(264, 8).Value = 0.05647491591813569
$ws.Cells.Item(264, 9).Value = 0.009594838179717874
$ws.Cells.Item(264, 10).Value = 5.647491591813569
$ws.Cells.Item(264, 11).Value = 0.9594838179717874
$ws.Cells.Item(264, 12).Value = 5.885968565631013

# Row 265
$ws.Cells.Item(265, 2).Value = 1
$c = $ws.Cells.Item(265, 3)
$c.NumberFormat = "@"
$c.Value = "14"
$c.NumberFormat = "General"
$ws.Cells.Item(265, 8).Value = 0.08316442454385897
$ws.Cells.Item(265, 9).Value = 0.00817659926446897
$ws.Cells.Item(265, 10).Value = 8.316442454385896
$ws.Cells.Item(265, 11).Value = 0.8176599264468971
$ws.Cells.Item(265, 12).Value = 10.17102854792531

# Row 266
$c = $ws.Cells.Item(266, 3)
$c.NumberFormat = "@"
$c.Value = "3"
$c.NumberFormat = "General"
$ws.Cells.Item(266, 8).Value = 0.06472084493618491
$ws.Cells.Item(266, 9).Value = 0.005292729966542757
$ws.Cells.Item(266, 10).Value = 6.472084493618491
$ws.Cells.Item(266, 11).Value = 0.5292729966542756
$ws.Cells.Item(266, 12).Value = 12.22825372639613

# Row 267
$c = $ws.Cells.Item(267, 3)
$c.NumberFormat = "@"
$c.Value = "3"
$c.NumberFormat = "General"
$ws.Cells.Item(267, 8).Value = 0.0303998196442723
$ws.Cells.Item(267, 9).Value = 0.005023277609478475
$ws.Cells.Item(267, 10).Value = 3.03998196442723
$ws.Cells.Item(267, 11).Value = 0.5023277609478475
$ws.Cells.Item(267, 12).Value = 6.051789689447097

# Row 268
$c = $ws.Cells.Item(268, 3)
$c.NumberFormat = "@"
$c.Value = "6"
$c.NumberFormat = "General"
$ws.Cells.Item(268, 8).Value = 0.04442946489350508
$ws.Cells.Item(268, 9).Value = 0.002145720523664109
$ws.Cells.Item(268, 10).Value = 4.442946489350508
$ws.Cells.Item(268, 11).Value = 0.2145720523664109
$ws.Cells.Item(268, 12).Value = 20.70608189813823

# Row 269
$c = $ws.Cells.Item(269, 3)
$c.NumberFormat = "@"
$c.Value = "3"
$c.NumberFormat = "General"
$ws.Cells.Item(269, 8).Value = 0.05252604205079781
$ws.Cells.Item(269, 9).Value = 0.005416085979678205
$ws.Cells.Item(269, 10).Value = 5.252604205079781
$ws.Cells.Item(269, 11).Value = 0.5416085979678205
$ws.Cells.Item(269, 12).Value = 9.698155134147008

# Row 270
$c = $ws.Cells.Item(270, 3)
$c.NumberFormat = "@"
$c.Value = "3"
$c.NumberFormat = "General"
$ws.Cells.Item(270, 8).Value = 0.02503470039176614
$ws.Cells.Item(270, 9).Value = 0.005075999754338674
$ws.Cells.Item(270, 10).Value = 2.503470039176614
$ws.Cells.Item(270, 11).Value = 0.5075999754338675
$ws.Cells.Item(270, 12).Value = 4.931974311142925

# Row 271
$c = $ws.Cells.Item(271, 3)
$c.NumberFormat = "@"
$c.Value = "6"
$c.NumberFormat = "General"
$ws.Cells.Item(271, 8).Value = 0.03775238396309422
$ws.Cells.Item(271, 9).Value = 0.002215644267868895
$ws.Cells.Item(271, 10).Value = 3.775238396309422
$ws.Cells.Item(271, 11).Value = 0.2215644267868895
$ws.Cells.Item(271, 12).Value = 17.03900960572797

# Row 272
$c = $ws.Cells.Item(272, 3)
$c.NumberFormat = "@"
$c.Value = "3"
$c.NumberFormat = "General"
$ws.Cells.Item(272, 8).Value = 0.09557450583412574
$ws.Cells.Item(272, 9).Value = 0.01471874494530004
$ws.Cells.Item(272, 10).Value = 9.557450583412574
$ws.Cells.Item(272, 11).Value = 1.471874494530004
$ws.Cells.Item(272, 12).Value = 6.493386915074195

# Row 273
$c = $ws.Cells.Item(273, 3)
$c.NumberFormat = "@"
$c.Value = "3"
$c.NumberFormat = "General"
$ws.Cells.Item(273, 8).Value = 0.05372824302963108
$ws.Cells.Item(273, 9).Value = 0.01050726168263165
$ws.Cells.Item(273, 10).Value = 5.372824302963108
$ws.Cells.Item(273, 11).Value = 1.050726168263165
$ws.Cells.Item(273, 12).Value = 5.113439129287424

# Row 274
$c = $ws.Cells.Item(274, 3)
$c.NumberFormat = "@"
$c.Value = "6"
$c.NumberFormat = "General"
$ws.Cells.Item(274, 8).Value = 0.05679768375088701
$ws.Cells.Item(274, 9).Value = 0.003666802204926504
$ws.Cells.Item(274, 10).Value = 5.679768375088701
$ws.Cells.Item(274, 11).Value = 0.3666802204926504
$ws.Cells.Item(274, 12).Value = 15.48970480997773

# Row 275
$c = $ws.Cells.Item(275, 3)
$c.NumberFormat = "@"
$c.Value = "3"
$c.NumberFormat = "General"
$ws.Cells.Item(275, 8).Value = 0.08830587204638318
$ws.Cells.Item(275, 9).Value = 0.002814341536194423
$ws.Cells.Item(275, 10).Value = 8.830587204638318
$ws.Cells.Item(275, 11).Value = 0.2814341536194423
$ws.Cells.Item(275, 12).Value = 31.37709866080829

# Row 276
$c = $ws.Cells.Item(276, 3)
$c.NumberFormat = "@"
$c.Value = "3"
$c.NumberFormat = "General"
$ws.Cells.Item(276, 8).Value = 0.05252604205079781
$ws.Cells.Item(276, 9).Value = 0.02046076925656211
$ws.Cells.Item(276, 10).Value = 5.252604205079781
$ws.Cells.Item(276, 11).Value = 2.046076925656211
$ws.Cells.Item(276, 12).Value = 2.56715871198009

# Row 277
$c = $ws.Cells.Item(277, 3)
$c.NumberFormat = "@"
$c.Value = "6"
$c.NumberFormat = "General"
$ws.Cells.Item(277, 8).Value = 0.05425359234869997
$ws.Cells.Item(277, 9).Value = 0.0080623909101685
$ws.Cells.Item(277, 10).Value = 5.425359234869997
$ws.Cells.Item(277, 11).Value = 0.8062390910168501
$ws.Cells.Item(277, 12).Value = 6.729218783013102

# Row 278
$c = $ws.Cells.Item(278, 3)
$c.NumberFormat = "@"
$c.Value = "3"
$c.NumberFormat = "General"
$ws.Cells.Item(278, 8).Value = 0.06882159843166558
$ws.Cells.Item(278, 9).Value = 0.002917885881604328
$ws.Cells.Item(278, 10).Value = 6.882159843166558
$ws.Cells.Item(278, 11).Value = 0.2917885881604328
$ws.Cells.Item(278, 12).Value = 23.5861172177939

# Row 279
$c = $ws.Cells.Item(279, 3)
$c.NumberFormat = "@"
$c.Value = "3"
$c.NumberFormat = "General"
$ws.Cells.Item(279, 8).Value = 0.04920574856226789
$ws.Cells.Item(279, 9).Value = 0.04148374821155287
$ws.Cells.Item(279, 10).Value = 4.920574856226789
$ws.Cells.Item(279, 11).Value = 4.148374821155286
$ws.Cells.Item(279, 12).Value = 1.186145193807837

# Row 280
$c = $ws.Cells.Item(280, 3)
$c.NumberFormat = "@"
$c.Value = "6"
$c.NumberFormat = "General"
$ws.Cells.Item(280, 8).Value = 0.04880884817015163
$ws.Cells.Item(280, 9).Value = 0.02009363307914265
$ws.Cells.Item(280, 10).Value = 4.880884817015163
$ws.Cells.Item(280, 11).Value = 2.009363307914265
$ws.Cells.Item(280, 12).Value = 2.429070341730068
